# fix(module3): use uncon_planned_qty for future production; keep produced for today
# Rewrites rows 4-6 (re-mapping material/location combos) and appends a new
# row 7 for MAT_B / PLANT_001 on the NetDemand sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: MAT_A / DC_001 - quantity changes
$ws.Range("F2").Value = -239

# Row 3: MAT_A / DC_002 - quantity changes
$ws.Range("F3").Value = -524

# Row 4: was MAT_B / DC_001 -> becomes MAT_A / PLANT_001
$ws.Range("A4").Value = "MAT_A"
$ws.Range("B4").Value = "PLANT_001"
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = -832
$ws.Range("H4").Value = 1

# Row 5: was MAT_B / DC_002 -> becomes MAT_B / DC_001
$ws.Range("B5").Value = "DC_001"
$ws.Range("F5").Value = -117
$ws.Range("H5").Value = 4

# Row 6: was MAT_B / PLANT_001 -> becomes MAT_B / DC_002
$ws.Range("B6").Value = "DC_002"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = -45

# Row 7 (new): MAT_B / PLANT_001
$ws.Range("A7").Value = "MAT_B"
$ws.Range("B7").Value = "PLANT_001"
$ws.Range("C7").Value = 45299
$ws.Range("C7").NumberFormat = $ws.Range("C6").NumberFormat
$ws.Range("D7").Value = "Distribution Demand - Forecast"
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = -117
$ws.Range("G7").Value = 45298
$ws.Range("G7").NumberFormat = $ws.Range("G6").NumberFormat
$ws.Range("H7").Value = 1
